# Apply the MegaSena data update: append lottery draws 2966-2971 and a
# trailing "spacer" row (429) with a single-space shared string in B:G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New draw rows (Concurso, Bola1..Bola6)
$rows = @(
    @(2966, 6, 7, 9, 43, 44, 53),
    @(2967, 1, 6, 38, 47, 56, 60),
    @(2968, 10, 11, 22, 26, 36, 46),
    @(2969, 1, 2, 5, 14, 18, 32),
    @(2970, 22, 32, 37, 41, 42, 59),
    @(2971, 1, 27, 39, 40, 46, 56)
)

$startRow = 423
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Trailing row with a single space in columns B:G (no value in A)
$lastRow = $startRow + $rows.Length
$ws.Range("B${lastRow}:G${lastRow}").Value = " "

# Keep the selection/active cell consistent with the appended data
$ws.Range("E" + ($lastRow + 2)).Select()
